# Auto-generated Excel COM-interop script
# Applies numeric corrections to the Leve profit-tracking sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the authoritative diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 6173462  # H43: was 18519018
$ws.Cells.Item(43, 9).Value = 800.2  # I43: was 750
$ws.Cells.Item(43, 10).Value = 13889289  # J43: was 55555556
$ws.Cells.Item(43, 11).Value = 800.2  # K43: was 750
$ws.Cells.Item(43, 12).Value = 13889289  # L43: was 55555556
$ws.Cells.Item(43, 13).Value = -731.2  # M43: was -681
$ws.Cells.Item(43, 14).Value = -13889427  # N43: was -55555694

$ws.Cells.Item(58, 8).Value = 3060  # H58: was 3160.5652
$ws.Cells.Item(58, 9).Value = 444.6  # I58: was 673.5
$ws.Cells.Item(58, 10).Value = 3994.0715  # J58: was 3684.158
$ws.Cells.Item(58, 11).Value = 1333.8  # K58: was 2020.5
$ws.Cells.Item(58, 12).Value = 11982.2145  # L58: was 11052.474
$ws.Cells.Item(58, 13).Value = -1183.8  # M58: was -1870.5
$ws.Cells.Item(58, 14).Value = -12282.2145  # N58: was -11352.474

$ws.Cells.Item(76, 8).Value = 3255  # H76: was 3433.2
$ws.Cells.Item(76, 9).Value = 3005.7144  # I76: was 3041.5
$ws.Cells.Item(76, 11).Value = 3005.7144  # K76: was 3041.5
$ws.Cells.Item(76, 13).Value = -2690.7144  # M76: was -2726.5

$ws.Cells.Item(79, 8).Value = 3255  # H79: was 3433.2
$ws.Cells.Item(79, 9).Value = 3005.7144  # I79: was 3041.5
$ws.Cells.Item(79, 11).Value = 3005.7144  # K79: was 3041.5
$ws.Cells.Item(79, 13).Value = -1913.7144  # M79: was -1949.5

$ws.Cells.Item(87, 8).Value = 44249.5  # H87: was 49996.668
$ws.Cells.Item(87, 10).Value = 44249.5  # J87: was 49996.668
$ws.Cells.Item(87, 12).Value = 44249.5  # L87: was 49996.668
$ws.Cells.Item(87, 14).Value = -46745.5  # N87: was -52492.668

$ws.Cells.Item(90, 8).Value = 44249.5  # H90: was 49996.668
$ws.Cells.Item(90, 10).Value = 44249.5  # J90: was 49996.668
$ws.Cells.Item(90, 12).Value = 132748.5  # L90: was 149990.004
$ws.Cells.Item(90, 14).Value = -145228.5  # N90: was -162470.004

$ws.Cells.Item(111, 8).Value = 2340.7334  # H111: was 3743.5293
$ws.Cells.Item(111, 9).Value = 2159  # I111: was 2017.625
$ws.Cells.Item(111, 10).Value = 2499.75  # J111: was 5277.6665
$ws.Cells.Item(111, 11).Value = 6477  # K111: was 6052.875
$ws.Cells.Item(111, 12).Value = 7499.25  # L111: was 15832.9995
$ws.Cells.Item(111, 13).Value = -3410  # M111: was -2985.875
$ws.Cells.Item(111, 14).Value = -13633.25  # N111: was -21966.9995

$ws.Cells.Item(137, 8).Value = 2235.9092  # H137: was 2221.5
$ws.Cells.Item(137, 9).Value = 1398.6  # I137: was 1314.8334
$ws.Cells.Item(137, 10).Value = 2933.6667  # J137: was 2901.5
$ws.Cells.Item(137, 11).Value = 4195.799999999999  # K137: was 3944.5002
$ws.Cells.Item(137, 12).Value = 8801.000100000001  # L137: was 8704.5
$ws.Cells.Item(137, 13).Value = -1645.799999999999  # M137: was -1394.5002
$ws.Cells.Item(137, 14).Value = -13901.0001  # N137: was -13804.5

$ws.Cells.Item(138, 8).Value = 2021.8586  # H138: was 2080.42
$ws.Cells.Item(138, 9).Value = 849.0769  # I138: was 994.5
$ws.Cells.Item(138, 10).Value = 2199.1396  # J138: was 2201.078
$ws.Cells.Item(138, 11).Value = 2547.2307  # K138: was 2983.5
$ws.Cells.Item(138, 12).Value = 6597.418799999999  # L138: was 6603.234
$ws.Cells.Item(138, 13).Value = 2592.7693  # M138: was 2156.5
$ws.Cells.Item(138, 14).Value = -16877.4188  # N138: was -16883.234

$ws.Cells.Item(141, 8).Value = 9235.615  # H141: was 9281.691999999999
$ws.Cells.Item(141, 9).Value = 10496.637  # I141: was 10551.091
$ws.Cells.Item(141, 11).Value = 31489.911  # K141: was 31653.273
$ws.Cells.Item(141, 13).Value = -26309.911  # M141: was -26473.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5052.4116  # H32: was 5781.963
$ws.Cells.Item(32, 9).Value = 5025.0938  # I32: was 5805.36
$ws.Cells.Item(32, 11).Value = 5025.0938  # K32: was 5805.36
$ws.Cells.Item(32, 13).Value = -4738.0938  # M32: was -5518.36

$ws.Cells.Item(132, 8).Value = 2112.6304  # H132: was 2285.4524
$ws.Cells.Item(132, 9).Value = 1814.2195  # I132: was 1978.1351
$ws.Cells.Item(132, 11).Value = 5442.6585  # K132: was 5934.4053
$ws.Cells.Item(132, 13).Value = -2912.6585  # M132: was -3404.4053

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1636.8  # H20: was 1696.5652
$ws.Cells.Item(20, 9).Value = 1383.4445  # I20: was 1437.6875
$ws.Cells.Item(20, 11).Value = 1383.4445  # K20: was 1437.6875
$ws.Cells.Item(20, 13).Value = -1136.4445  # M20: was -1190.6875

$ws.Cells.Item(86, 8).Value = 3114.8125  # H86: was 3181.9688
$ws.Cells.Item(86, 9).Value = 3681.9  # I86: was 3789.35
$ws.Cells.Item(86, 11).Value = 3681.9  # K86: was 3789.35
$ws.Cells.Item(86, 13).Value = -2558.9  # M86: was -2666.35

$ws.Cells.Item(87, 8).Value = 75000  # H87: was 59666.668
$ws.Cells.Item(87, 10).Value = 75000  # J87: was 59666.668
$ws.Cells.Item(87, 12).Value = 75000  # L87: was 59666.668
$ws.Cells.Item(87, 14).Value = -77496  # N87: was -62162.668

$ws.Cells.Item(89, 8).Value = 3114.8125  # H89: was 3181.9688
$ws.Cells.Item(89, 9).Value = 3681.9  # I89: was 3789.35
$ws.Cells.Item(89, 11).Value = 18409.5  # K89: was 18946.75
$ws.Cells.Item(89, 13).Value = -12793.5  # M89: was -13330.75

$ws.Cells.Item(90, 8).Value = 75000  # H90: was 59666.668
$ws.Cells.Item(90, 10).Value = 75000  # J90: was 59666.668
$ws.Cells.Item(90, 12).Value = 225000  # L90: was 179000.004
$ws.Cells.Item(90, 14).Value = -237480  # N90: was -191480.004

$ws.Cells.Item(134, 8).Value = 5079.6333  # H134: was 8513.275
$ws.Cells.Item(134, 9).Value = 1788.64  # I134: was 5608.48
$ws.Cells.Item(134, 10).Value = 21534.6  # J134: was 26668.25
$ws.Cells.Item(134, 11).Value = 5365.92  # K134: was 16825.44
$ws.Cells.Item(134, 12).Value = 64603.8  # L134: was 80004.75
$ws.Cells.Item(134, 13).Value = -2830.92  # M134: was -14290.44
$ws.Cells.Item(134, 14).Value = -69673.79999999999  # N134: was -85074.75

$ws.Cells.Item(140, 8).Value = 22322.5  # H140: was 20709.092
$ws.Cells.Item(140, 10).Value = 22379.629  # J140: was 20702
$ws.Cells.Item(140, 12).Value = 22379.629  # L140: was 20702
$ws.Cells.Item(140, 14).Value = -32739.629  # N140: was -31062

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 1888.5  # H4: was 1147.5

$ws.Cells.Item(31, 8).Value = 1022.9655  # H31: was 1168.4814
$ws.Cells.Item(31, 9).Value = 796.8333  # I31: was 901.625
$ws.Cells.Item(31, 10).Value = 1393  # J31: was 1556.6364
$ws.Cells.Item(31, 11).Value = 796.8333  # K31: was 901.625
$ws.Cells.Item(31, 12).Value = 1393  # L31: was 1556.6364
$ws.Cells.Item(31, 13).Value = -501.8333  # M31: was -606.625
$ws.Cells.Item(31, 14).Value = -1983  # N31: was -2146.6364

$ws.Cells.Item(34, 8).Value = 1022.9655  # H34: was 1168.4814
$ws.Cells.Item(34, 9).Value = 796.8333  # I34: was 901.625
$ws.Cells.Item(34, 10).Value = 1393  # J34: was 1556.6364
$ws.Cells.Item(34, 11).Value = 796.8333  # K34: was 901.625
$ws.Cells.Item(34, 12).Value = 1393  # L34: was 1556.6364
$ws.Cells.Item(34, 13).Value = -594.8333  # M34: was -699.625
$ws.Cells.Item(34, 14).Value = -1797  # N34: was -1960.6364

$ws.Cells.Item(58, 8).Value = 1174.5927  # H58: was 1252.5238
$ws.Cells.Item(58, 9).Value = 1153.3684  # I58: was 1235.9286
$ws.Cells.Item(58, 10).Value = 1225  # J58: was 1285.7142
$ws.Cells.Item(58, 11).Value = 1153.3684  # K58: was 1235.9286
$ws.Cells.Item(58, 12).Value = 1225  # L58: was 1285.7142
$ws.Cells.Item(58, 13).Value = -950.3684000000001  # M58: was -1032.9286
$ws.Cells.Item(58, 14).Value = -1631  # N58: was -1691.7142

$ws.Cells.Item(105, 8).Value = 627.1875  # H105: was 671
$ws.Cells.Item(105, 9).Value = 563.25  # I105: was 617.1818
$ws.Cells.Item(105, 11).Value = 563.25  # K105: was 617.1818
$ws.Cells.Item(105, 13).Value = 1183.75  # M105: was 1129.8182

$ws.Cells.Item(132, 8).Value = 8912.556  # H132: was 7601
$ws.Cells.Item(132, 9).Value = 12146.728  # I132: was 8663.0625
$ws.Cells.Item(132, 10).Value = 3830.2856  # J132: was 4202.4
$ws.Cells.Item(132, 11).Value = 36440.18399999999  # K132: was 25989.1875
$ws.Cells.Item(132, 12).Value = 11490.8568  # L132: was 12607.2
$ws.Cells.Item(132, 13).Value = -33910.18399999999  # M132: was -23459.1875
$ws.Cells.Item(132, 14).Value = -16550.8568  # N132: was -17667.2

$ws.Cells.Item(134, 8).Value = 9525148  # H134: was 10754111
$ws.Cells.Item(134, 9).Value = 11906034  # I134: was 13334682
$ws.Cells.Item(134, 10).Value = 1601  # J134: was 1732.3334
$ws.Cells.Item(134, 11).Value = 35718102  # K134: was 40004046
$ws.Cells.Item(134, 12).Value = 4803  # L134: was 5197.0002
$ws.Cells.Item(134, 13).Value = -35715567  # M134: was -40001511
$ws.Cells.Item(134, 14).Value = -9873  # N134: was -10267.0002

$ws.Cells.Item(136, 8).Value = 1174.5927  # H136: was 1252.5238
$ws.Cells.Item(136, 9).Value = 1153.3684  # I136: was 1235.9286
$ws.Cells.Item(136, 10).Value = 1225  # J136: was 1285.7142
$ws.Cells.Item(136, 11).Value = 3460.1052  # K136: was 3707.7858
$ws.Cells.Item(136, 12).Value = 3675  # L136: was 3857.1426
$ws.Cells.Item(136, 13).Value = -910.1052  # M136: was -1157.7858
$ws.Cells.Item(136, 14).Value = -8775  # N136: was -8957.142599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 164613.08  # H4: was 151633.5
$ws.Cells.Item(4, 9).Value = 172778.1  # I4: was 126757.266
$ws.Cells.Item(4, 10).Value = 160870.8  # J4: was 167857.12
$ws.Cells.Item(4, 11).Value = 518334.3  # K4: was 380271.798
$ws.Cells.Item(4, 12).Value = 482612.4  # L4: was 503571.36
$ws.Cells.Item(4, 13).Value = -518222.3  # M4: was -380159.798
$ws.Cells.Item(4, 14).Value = -482836.4  # N4: was -503795.36

$ws.Cells.Item(34, 8).Value = 1431.12  # H34: was 1599
$ws.Cells.Item(34, 9).Value = 614.2857  # I34: was 650
$ws.Cells.Item(34, 10).Value = 1748.7778  # J34: was 1822.2941
$ws.Cells.Item(34, 11).Value = 1842.8571  # K34: was 1950
$ws.Cells.Item(34, 12).Value = 5246.3334  # L34: was 5466.8823
$ws.Cells.Item(34, 13).Value = -1758.8571  # M34: was -1866
$ws.Cells.Item(34, 14).Value = -5414.3334  # N34: was -5634.8823

$ws.Cells.Item(56, 8).Value = 6807.6924  # H56: was 7265.643
$ws.Cells.Item(56, 9).Value = 6807.6924  # I56: was 7265.643
$ws.Cells.Item(56, 11).Value = 6807.6924  # K56: was 7265.643
$ws.Cells.Item(56, 13).Value = -6277.6924  # M56: was -6735.643

$ws.Cells.Item(61, 8).Value = 300  # H61: was 252.6
$ws.Cells.Item(61, 9).Value = 0  # I61: was 98
$ws.Cells.Item(61, 10).Value = 300  # J61: was 291.25
$ws.Cells.Item(61, 11).Value = 0  # K61: was 294
$ws.Cells.Item(61, 12).Value = 900  # L61: was 873.75
$ws.Cells.Item(61, 13).Value = ""  # M61: clear (was -79)
$ws.Cells.Item(61, 14).Value = -1330  # N61: was -1303.75

$ws.Cells.Item(68, 8).Value = 2202  # H68: was 2237.4285
$ws.Cells.Item(68, 9).Value = 855.75  # I68: was 874.5
$ws.Cells.Item(68, 10).Value = 3279  # J68: was 3076.1538
$ws.Cells.Item(68, 11).Value = 2567.25  # K68: was 2623.5
$ws.Cells.Item(68, 12).Value = 9837  # L68: was 9228.4614
$ws.Cells.Item(68, 13).Value = -1756.25  # M68: was -1812.5
$ws.Cells.Item(68, 14).Value = -11459  # N68: was -10850.4614

$ws.Cells.Item(71, 8).Value = 2202  # H71: was 2237.4285
$ws.Cells.Item(71, 9).Value = 855.75  # I71: was 874.5
$ws.Cells.Item(71, 10).Value = 3279  # J71: was 3076.1538
$ws.Cells.Item(71, 11).Value = 7701.75  # K71: was 7870.5
$ws.Cells.Item(71, 12).Value = 29511  # L71: was 27685.3842
$ws.Cells.Item(71, 13).Value = -3645.75  # M71: was -3814.5
$ws.Cells.Item(71, 14).Value = -37623  # N71: was -35797.3842

$ws.Cells.Item(107, 8).Value = 3515.4  # H107: was 3609.9707
$ws.Cells.Item(107, 9).Value = 567.26086  # I107: was 579.4091
$ws.Cells.Item(107, 11).Value = 1701.78258  # K107: was 1738.2273
$ws.Cells.Item(107, 13).Value = 218.2174199999999  # M107: was 181.7727

$ws.Cells.Item(129, 8).Value = 18117110  # H129: was 21931084
$ws.Cells.Item(129, 9).Value = 33334060  # I129: was 37037796
$ws.Cells.Item(129, 10).Value = 6411766  # J129: was 8335043.5
$ws.Cells.Item(129, 11).Value = 100002180  # K129: was 111113388
$ws.Cells.Item(129, 12).Value = 19235298  # L129: was 25005130.5
$ws.Cells.Item(129, 13).Value = -99997180  # M129: was -111108388
$ws.Cells.Item(129, 14).Value = -19245298  # N129: was -25015130.5

$ws.Cells.Item(131, 8).Value = 23811036  # H131: was 24391776
$ws.Cells.Item(131, 9).Value = 166667280  # I131: was 200000600
$ws.Cells.Item(131, 11).Value = 500001840  # K131: was 600001800
$ws.Cells.Item(131, 13).Value = -499996800  # M131: was -599996760

$ws.Cells.Item(138, 8).Value = 2868.5557  # H138: was 2101.1333
$ws.Cells.Item(138, 9).Value = 2997.7144  # I138: was 2073.6667
$ws.Cells.Item(138, 10).Value = 2416.5  # J138: was 2211
$ws.Cells.Item(138, 11).Value = 8993.143199999999  # K138: was 6221.000100000001
$ws.Cells.Item(138, 12).Value = 7249.5  # L138: was 6633
$ws.Cells.Item(138, 13).Value = -3853.143199999999  # M138: was -1081.000100000001
$ws.Cells.Item(138, 14).Value = -17529.5  # N138: was -16913

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 0  # H5: was 2000
$ws.Cells.Item(5, 9).Value = 0  # I5: was 2000
$ws.Cells.Item(5, 11).Value = 0  # K5: was 2000
$ws.Cells.Item(5, 13).Value = ""  # M5: clear (was -1888)

$ws.Cells.Item(52, 8).Value = 19995  # H52: was 20000
$ws.Cells.Item(52, 10).Value = 19995  # J52: was 20000
$ws.Cells.Item(52, 12).Value = 19995  # L52: was 20000
$ws.Cells.Item(52, 14).Value = -20513  # N52: was -20518

$ws.Cells.Item(70, 8).Value = 40913756  # H70: was 75005384
$ws.Cells.Item(70, 9).Value = 41670900  # I70: was 125003000
$ws.Cells.Item(70, 10).Value = 40005180  # J70: was 50006576
$ws.Cells.Item(70, 11).Value = 41670900  # K70: was 125003000
$ws.Cells.Item(70, 12).Value = 40005180  # L70: was 50006576
$ws.Cells.Item(70, 13).Value = -41670630  # M70: was -125002730
$ws.Cells.Item(70, 14).Value = -40005720  # N70: was -50007116

$ws.Cells.Item(73, 8).Value = 40913756  # H73: was 75005384
$ws.Cells.Item(73, 9).Value = 41670900  # I73: was 125003000
$ws.Cells.Item(73, 10).Value = 40005180  # J73: was 50006576
$ws.Cells.Item(73, 11).Value = 41670900  # K73: was 125003000
$ws.Cells.Item(73, 12).Value = 40005180  # L73: was 50006576
$ws.Cells.Item(73, 13).Value = -41669964  # M73: was -125002064
$ws.Cells.Item(73, 14).Value = -40007052  # N73: was -50008448

$ws.Cells.Item(132, 8).Value = 2164.5806  # H132: was 2325.6785
$ws.Cells.Item(132, 9).Value = 1719.238  # I132: was 1832.6842
$ws.Cells.Item(132, 10).Value = 3099.8  # J132: was 3366.4443
$ws.Cells.Item(132, 11).Value = 5157.714  # K132: was 5498.0526
$ws.Cells.Item(132, 12).Value = 9299.400000000001  # L132: was 10099.3329
$ws.Cells.Item(132, 13).Value = -2627.714  # M132: was -2968.0526
$ws.Cells.Item(132, 14).Value = -14359.4  # N132: was -15159.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 2805.9092  # H2: was 2803.5
$ws.Cells.Item(2, 9).Value = 1980  # I2: was 3980
$ws.Cells.Item(2, 10).Value = 2888.5  # J2: was 2696.5454
$ws.Cells.Item(2, 11).Value = 1980  # K2: was 3980
$ws.Cells.Item(2, 12).Value = 2888.5  # L2: was 2696.5454
$ws.Cells.Item(2, 13).Value = -1868  # M2: was -3868
$ws.Cells.Item(2, 14).Value = -3112.5  # N2: was -2920.5454

$ws.Cells.Item(122, 8).Value = 10901195  # H122: was 11809341
$ws.Cells.Item(122, 9).Value = 18892408  # I122: was 21798404
$ws.Cells.Item(122, 11).Value = 56677224  # K122: was 65395212
$ws.Cells.Item(122, 13).Value = -56674774  # M122: was -65392762

$ws.Cells.Item(132, 8).Value = 50033.383  # H132: was 58117.668
$ws.Cells.Item(132, 9).Value = 2153.9092  # I132: was 2311
$ws.Cells.Item(132, 10).Value = 102700.8  # J132: was 127876
$ws.Cells.Item(132, 11).Value = 6461.7276  # K132: was 6933
$ws.Cells.Item(132, 12).Value = 308102.4  # L132: was 383628
$ws.Cells.Item(132, 13).Value = -3931.7276  # M132: was -4403
$ws.Cells.Item(132, 14).Value = -313162.4  # N132: was -388688

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 3500  # H2: was 5000
$ws.Cells.Item(2, 10).Value = 2000  # J2: was 0
$ws.Cells.Item(2, 12).Value = 2000  # L2: was 0
$ws.Cells.Item(2, 14).Value = -2224  # N2: new cell (added)

$ws.Cells.Item(128, 8).Value = 63330  # H128: was 74994.5
$ws.Cells.Item(128, 10).Value = 63330  # J128: was 74994.5
$ws.Cells.Item(128, 12).Value = 63330  # L128: was 74994.5
$ws.Cells.Item(128, 14).Value = -73290  # N128: was -84954.5

$ws.Cells.Item(132, 8).Value = 4362.64  # H132: was 4695.0435
$ws.Cells.Item(132, 9).Value = 5510.25  # I132: was 6220.2856
$ws.Cells.Item(132, 11).Value = 16530.75  # K132: was 18660.8568
$ws.Cells.Item(132, 13).Value = -14000.75  # M132: was -16130.8568
